$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (G) values, replacing the former Strike# values.
$values = @{
    2  = 3
    3  = 2
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 3
    10 = 0
    11 = 0
    12 = 1
    13 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
